# Apply the "TCDs" fixes described in the commit:
#  - Update the rule-15-5d description text (R5) to the new, corrected wording
#  - Fix several transformer kVA rating test values in row 13 and row 21
#    (they were off by a factor of 1000, e.g. 25 -> 25000)
#  - Correct the expected efficiency value in Q38 (0.983 -> 0.9823)
#  - Move the active selection to Q9 on the TCDs sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TCDs")

# --- Row 5: corrected rule description text -------------------------------
$ws.Range("R5").Value = "Dry Three-Phase transformer modeled in the User Model whose efficiency is > than Table 8.4.4. Baseline RMR does not match Table 8.4.4"

# --- Row 13: kVA values corrected (x1000) ----------------------------------
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 25000
$ws.Range("Q13").Value = 30000
$ws.Range("R13").Value = 30000
$ws.Range("W13").Value = 30000
$ws.Range("X13").Value = 30000
$ws.Range("Y13").Value = 25000
$ws.Range("Z13").Value = 25000
$ws.Range("AA13").Value = 30000
$ws.Range("AB13").Value = 30000
$ws.Range("AC13").Value = 25000
$ws.Range("AD13").Value = 25000
$ws.Range("AE13").Value = 25000
$ws.Range("AF13").Value = 25000
$ws.Range("AG13").Value = 25000

# --- Row 21: kVA values corrected (x1000) ----------------------------------
$ws.Range("S21").Value = 25000
$ws.Range("T21").Value = 25000
$ws.Range("U21").Value = 25000
$ws.Range("V21").Value = 25000

# --- Row 38: corrected efficiency value ------------------------------------
$ws.Range("Q38").Value = 0.9823

# --- Update the selected / active cell on the TCDs sheet -------------------
$ws.Activate()
$ws.Range("Q9").Select()
